# TestData.xlsx edit:
#   - "LoginTestData" sheet renamed to "Sheet1"
#   - extra duplicate data rows (5-7) removed, shrinking the used range
#     from A1:C7 down to A1:C4 (rows shift up automatically)
#   - the hyperlinks that pointed at the now-removed rows (A5:A7) are
#     cleaned up, leaving only the hyperlinks on A2:A4
#   - selection moved to J13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet.
$ws.Name = "Sheet1"

# Remove the extra rows (5,6,7); remaining rows below shift up to fill
# the gap, same as selecting the rows in Excel and choosing Delete.
$ws.Range("A5:C7").Delete()

# The hyperlink collection keeps stale entries for the rows that no
# longer exist (A5:A7) - drop exactly those, leaving A2:A4 untouched so
# their formatting/style is not disturbed. Re-query the live collection
# on every deletion since removing an item shifts the others.
$maxIterations = $ws.Hyperlinks.Count + 1
for ($i = 0; $i -lt $maxIterations; $i++) {
  $staleLink = $null
  foreach ($link in $ws.Hyperlinks) {
    if ($link.Range.Row -ge 5) {
      $staleLink = $link
      break
    }
  }
  if ($staleLink -eq $null) {
    break
  }
  $staleLink.Delete()
}

# Move the active selection.
[void]$ws.Range("J13").Select()
